$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the first item row: it used to be "sample 1 " (id 1, price 10),
# now it's "RTX 4090" (id 0, price 110000). Description stays "chu 1".
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "RTX 4090"
$ws.Range("C2").Value = 110000

# Every subsequent row's Item id shifts down by one (2->1, 3->2, ... 20->19).
for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1
}

# Move the active selection to C3.
$ws.Range("C3").Select()
